$d = $word.ActiveDocument

# The edit is authored by Gordon Chalmers and recorded as tracked-change
# insertions (w:ins) -- turn on Track Changes and set the author identity
# before making any edits.
$word.Application.UserName = "Gordon Chalmers"
$d.TrackRevisions = $true

# Locate the end of "...mathematical techniques." (the last sentence of
# the paragraph right before the one that begins "Due to my background
# and multiple PhDs ...") and append four new tracked-insertion runs
# there, reproducing the four separate w:ins blocks from the diff.
$full = $d.Content.Text
$anchor = "mathematical techniques."
$idx = $full.IndexOf($anchor)
if ($idx -lt 0) {
    throw "anchor text not found"
}
$endPos = $idx + $anchor.Length

# Each chunk becomes its own tracked insertion (separate w:ins wrapper),
# exactly like the four consecutive <w:ins> runs added in the diff.
$chunks = @(
    "  Al",
    "l of the graduate students I have worked ",
    "with ",
    "(except Xu Yang who is still in school) have graduated and gone on to successful careers in molecular work."
)

foreach ($chunk in $chunks) {
    $rng = $d.Range($endPos, $endPos)
    $rng.Select()
    $word.Selection.InsertAfter($chunk)
    $endPos = $endPos + $chunk.Length
}

Write-Output "inserted mentoring sentence"
